$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Datos actualizados" timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 15:52"

# --- Rows 4 (Estados Unidos): refresh totals ---
$ws.Cells.Item(4, 2).Value = 1064836
$ws.Cells.Item(4, 3).Value = 642
$ws.Cells.Item(4, 5).Value = 855683

# --- Row 9 (Alemania): refresh totals ---
$ws.Cells.Item(9, 2).Value = 161845
$ws.Cells.Item(9, 3).Value = 306
$ws.Cells.Item(9, 5).Value = 31875
$ws.Cells.Item(9, 7).Value = 3
$ws.Cells.Item(9, 8).Value = 6470

# --- Row 42 (Serbia): refresh "Muertes hoy" ---
$ws.Cells.Item(42, 6).Value = 71

# --- Row 61 (Kazajistan): refresh totals ---
$ws.Cells.Item(61, 4).Value = 858
$ws.Cells.Item(61, 5).Value = 2473

# --- Azerbaiyan overtakes Islandia: swap the two countries (rows 74-75) ---
# and refresh Azerbaiyan's case numbers.
$ws.Cells.Item(74, 1).Value = "Azerbaiyan"
$ws.Cells.Item(74, 2).Value = 1804
$ws.Cells.Item(74, 3).Value = 38
$ws.Cells.Item(74, 4).Value = 1325
$ws.Cells.Item(74, 5).Value = 455
$ws.Cells.Item(74, 6).Value = 17
$ws.Cells.Item(74, 7).Value = 1
$ws.Cells.Item(74, 8).Value = 24

$ws.Cells.Item(75, 1).Value = "Islandia"
$ws.Cells.Item(75, 2).Value = 1797
$ws.Cells.Item(75, 3).Value = 0
$ws.Cells.Item(75, 4).Value = 1656
$ws.Cells.Item(75, 5).Value = 131
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 0
$ws.Cells.Item(75, 8).Value = 10

# --- Suazilandia overtakes Aruba, Zambia and Monaco: insert it right
# after Bermudas (row 150) with refreshed numbers, and push Aruba,
# Zambia and Monaco down one row each (rows 151-154). Liechtenstein
# (row 155) is unaffected. ---
$ws.Cells.Item(151, 1).Value = "Suazilandia"
$ws.Cells.Item(151, 2).Value = 100
$ws.Cells.Item(151, 3).Value = 9
$ws.Cells.Item(151, 4).Value = 10
$ws.Cells.Item(151, 5).Value = 89
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 7).Value = 0
$ws.Cells.Item(151, 8).Value = 1

$ws.Cells.Item(152, 1).Value = "Aruba"
$ws.Cells.Item(152, 2).Value = 100
$ws.Cells.Item(152, 3).Value = 0
$ws.Cells.Item(152, 4).Value = 73
$ws.Cells.Item(152, 5).Value = 25
$ws.Cells.Item(152, 6).Value = 4
$ws.Cells.Item(152, 7).Value = 0
$ws.Cells.Item(152, 8).Value = 2

$ws.Cells.Item(153, 1).Value = "Zambia"
$ws.Cells.Item(153, 2).Value = 97
$ws.Cells.Item(153, 3).Value = 0
$ws.Cells.Item(153, 4).Value = 54
$ws.Cells.Item(153, 5).Value = 40
$ws.Cells.Item(153, 6).Value = 1
$ws.Cells.Item(153, 7).Value = 0
$ws.Cells.Item(153, 8).Value = 3

$ws.Cells.Item(154, 1).Value = "Monaco"
$ws.Cells.Item(154, 2).Value = 95
$ws.Cells.Item(154, 3).Value = 0
$ws.Cells.Item(154, 4).Value = 58
$ws.Cells.Item(154, 5).Value = 33
$ws.Cells.Item(154, 6).Value = 1
$ws.Cells.Item(154, 7).Value = 0
$ws.Cells.Item(154, 8).Value = 4
